# Insert a new data row at row 244 (pushing the existing rows 244:341 down
# to 245:342), then populate the new row with the Ciruela / Angeleno entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("244:244").Insert()

$ws.Range("A244").Value = 10
$ws.Range("B244").Value = "Vega Modelo de Temuco"
$ws.Range("C244").Value = "La Araucanía"
$ws.Range("D244").Value = 45009
$ws.Range("E244").Value = 9
$ws.Range("F244").Value = "Fruta"
$ws.Range("G244").Value = 100103
$ws.Range("H244").Value = "Frutos de hueso (carozo)"
$ws.Range("I244").Value = 100103002
$ws.Range("J244").Value = "Ciruela"
$ws.Range("K244").Value = "Angeleno"
$ws.Range("L244").Value = "Primera"
$ws.Range("M244").Value = 200
$ws.Range("N244").Value = 13000
$ws.Range("O244").Value = 14000
$ws.Range("P244").Value = 13500
$ws.Range("Q244").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R244").Value = "Región de O'Higgins"
$ws.Range("S244").Value = 750
$ws.Range("T244").Value = 18
